$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-14 19:18:45'
$ws.Range('I2').Value = '33.9 mm'
$ws.Range('E3').Value = '2026-02-14 19:18:48'
$ws.Range('I3').Value = '13.9 mm'
$ws.Range('N3').Value = '-6.6 °C 18:59 TU'
$ws.Range('E4').Value = '2026-02-14 19:18:50'
$ws.Range('J4').Value = '996.2 hPa'
$ws.Range('O4').Value = '11.0 °C'
$ws.Range('E5').Value = '2026-02-14 19:18:53'
$ws.Range('I5').Value = '21.0 mm'
$ws.Range('N5').Value = '-6.3 °C 18:58 TU'
$ws.Range('O5').Value = '-5.1 °C'
$ws.Range('E6').Value = '2026-02-14 19:18:56'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '78%'
$ws.Range('J6').Value = '996.2 hPa'
$ws.Range('E7').Value = '2026-02-14 19:18:58'
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = '54%'
$ws.Range('J7').Value = '996.4 hPa'
$ws.Range('O7').Value = '13.2 °C'
$ws.Range('E8').Value = '2026-02-14 19:19:01'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '64%'
$ws.Range('J8').Value = '996.2 hPa'
$ws.Range('E9').Value = '2026-02-14 19:19:03'
$ws.Range('N9').Value = '10.6 °C 18:59 TU'
$ws.Range('O9').Value = '11.9 °C'
$ws.Range('E10').Value = '2026-02-14 19:19:06'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '79%'
$ws.Range('O10').Value = '10.1 °C'
$ws.Range('E11').Value = '2026-02-14 19:19:09'
$ws.Range('E12').Value = '2026-02-14 19:19:11'
$ws.Range('E13').Value = '2026-02-14 19:19:14'
$ws.Range('J13').Value = '998.8 hPa'
$ws.Range('E14').Value = '2026-02-14 19:19:16'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '51%'
$ws.Range('E15').Value = '2026-02-14 19:19:19'
$ws.Range('N15').Value = '10.3 °C 18:51 TU'
$ws.Range('E16').Value = '2026-02-14 19:19:22'
$ws.Range('N16').Value = '-8.4 °C 18:58 TU'
$ws.Range('O16').Value = '-5.9 °C'
$ws.Range('E17').Value = '2026-02-14 19:19:24'
$ws.Range('L17').Value = '55.1 km/h - 275º 18:36 TU'
$ws.Range('E18').Value = '2026-02-14 19:19:27'
$ws.Range('J18').Value = '996.4 hPa'
$ws.Range('E19').Value = '2026-02-14 19:19:30'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '80%'
$ws.Range('E20').Value = '2026-02-14 19:19:33'
$ws.Range('I20').Value = '3.5 mm'
$ws.Range('L20').Value = '111.6 km/h - 342º 18:32 TU'
$ws.Range('N20').Value = '-7.2 °C 18:48 TU'
$ws.Range('O20').Value = '-5.2 °C'
$ws.Range('E21').Value = '2026-02-14 19:19:36'
$ws.Range('I21').Value = '0.1 mm'
$ws.Range('J21').Value = '998.7 hPa'
$ws.Range('O21').Value = '5.4 °C'
$ws.Range('E22').Value = '2026-02-14 19:19:38'
$ws.Range('N22').Value = '-8.6 °C 18:34 TU'
$ws.Range('O22').Value = '-6.7 °C'
$ws.Range('E23').Value = '2026-02-14 19:19:41'
$ws.Range('I23').Value = '37.5 mm'
$ws.Range('N23').Value = '-8.0 °C 18:51 TU'
$ws.Range('O23').Value = '-5.9 °C'
$ws.Range('E24').Value = '2026-02-14 19:19:44'
$ws.Range('J24').Value = '1000.4 hPa'
$ws.Range('E25').Value = '2026-02-14 19:19:47'
$ws.Range('I25').Value = '12.9 mm'
$ws.Range('N25').Value = '-6.7 °C 18:59 TU'
$ws.Range('E26').Value = '2026-02-14 19:19:49'
$ws.Range('E27').Value = '2026-02-14 19:19:52'
$ws.Range('N27').Value = '-4.6 °C 18:54 TU'
$ws.Range('O27').Value = '-2.9 °C'
$ws.Range('E28').Value = '2026-02-14 19:19:55'
$ws.Range('H28').NumberFormat = '@'
$ws.Range('H28').Value = '70%'
$ws.Range('J28').Value = '996.2 hPa'
$ws.Range('L28').Value = '47.5 km/h - 311º 18:50 TU'
$ws.Range('O28').Value = '9.1 °C'
$ws.Range('E29').Value = '2026-02-14 19:19:58'
$ws.Range('O29').Value = '11.6 °C'
$ws.Range('E30').Value = '2026-02-14 19:20:01'
$ws.Range('J30').Value = '996.1 hPa'
$ws.Range('E31').Value = '2026-02-14 19:20:04'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '69%'
$ws.Range('J31').Value = '995.2 hPa'
$ws.Range('E32').Value = '2026-02-14 19:20:06'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '85%'
$ws.Range('E33').Value = '2026-02-14 19:20:09'
$ws.Range('J33').Value = '998.2 hPa'
$ws.Range('N33').Value = '1.3 °C 18:32 TU'
$ws.Range('O33').Value = '4.1 °C'
$ws.Range('E34').Value = '2026-02-14 19:20:12'
$ws.Range('N34').Value = '-4.1 °C 18:47 TU'
$ws.Range('O34').Value = '-2.1 °C'
$ws.Range('E35').Value = '2026-02-14 19:20:15'
$ws.Range('H35').NumberFormat = '@'
$ws.Range('H35').Value = '84%'
$ws.Range('J35').Value = '1003.1 hPa'
$ws.Range('E36').Value = '2026-02-14 19:20:17'
$ws.Range('J36').Value = '996.9 hPa'
$ws.Range('E37').Value = '2026-02-14 19:20:20'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '67%'
$ws.Range('J37').Value = '997.1 hPa'
$ws.Range('E38').Value = '2026-02-14 19:20:23'
$ws.Range('K38').Value = '12.0 MJ/m2'
$ws.Range('O38').Value = '10.3 °C'
$ws.Range('E39').Value = '2026-02-14 19:20:26'
$ws.Range('I39').Value = '12.1 mm'
$ws.Range('N39').Value = '-8.0 °C 18:59 TU'
$ws.Range('O39').Value = '-5.6 °C'
$ws.Range('E40').Value = '2026-02-14 19:20:29'
$ws.Range('I40').Value = '0.6 mm'
$ws.Range('J40').Value = '999.3 hPa'
$ws.Range('O40').Value = '7.2 °C'
$ws.Range('E41').Value = '2026-02-14 19:20:32'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '48%'
$ws.Range('J41').Value = '998.2 hPa'
$ws.Range('K41').Value = '14.1 MJ/m2'
$ws.Range('E42').Value = '2026-02-14 19:20:34'
$ws.Range('H42').NumberFormat = '@'
$ws.Range('H42').Value = '63%'
$ws.Range('E43').Value = '2026-02-14 19:20:37'
$ws.Range('E44').Value = '2026-02-14 19:20:40'
$ws.Range('I44').Value = '35.6 mm'
$ws.Range('N44').Value = '-6.5 °C 18:59 TU'
$ws.Range('O44').Value = '-5.3 °C'
$ws.Range('E45').Value = '2026-02-14 19:20:43'
$ws.Range('I45').Value = '13.4 mm'
$ws.Range('J45').Value = '1005.3 hPa'
$ws.Range('O45').Value = '3.0 °C'
$ws.Range('E46').Value = '2026-02-14 19:20:45'
